# Sample Campaign.xlsx - "Importing campaign spreadsheet works"
#
# Fills in a couple of player usernames that were missing, fixes up a
# couple of character races, and leaves the workbook with the NPC sheet
# active (the last thing touched while testing the import).

$wb = $excel.ActiveWorkbook

# --- Major Characters sheet -------------------------------------------------
$major = $wb.Worksheets.Item("Major Characters")
$major.Activate()

# Leox and Ganamede were missing their PlayerUsername values.
$major.Range("F2").Value = "abcdefghijklmnopqrstuvwxyz"
$major.Range("F5").Value = "abcdefghijklmnopqrstuvwxyz"

# Karg was miscategorized as a Dwarf - he is actually a Half-Orc.
$major.Range("B7").Value = "Half-Orc"

[void]$major.Range("B7").Select()

# --- NPC sheet ---------------------------------------------------------------
$npc = $wb.Worksheets.Item("NPC")
$npc.Activate()

# Rhea is specifically a Protector Aasimar.
$npc.Range("B4").Value = "Protector Aasimar"

[void]$npc.Range("B4").Select()
